$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("O2").Value = 1.15
$ws.Range("P2").Value = 2.8
$ws.Range("Q2").Value = 1.47
$ws.Range("S2").Value = 2.14
$ws.Range("T2").Value = 1.53
$ws.Range("U2").Value = 2.58
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 28
$ws.Range("AB2").Value = 15.5
$ws.Range("AD2").Value = 20
$ws.Range("AF2").Value = 15.5
$ws.Range("AK2").Value = 17
$ws.Range("AL2").Value = 28
$ws.Range("AN2").Value = 7

# Row 3
$ws.Range("G3").Value = 2.02
$ws.Range("J3").Value = 3.85
$ws.Range("K3").Value = 4.3
$ws.Range("L3").Value = 1.27
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 4.6
$ws.Range("O3").Value = 1.23
$ws.Range("P3").Value = 2.24
$ws.Range("R3").Value = 1.5
$ws.Range("U3").Value = 2.3
$ws.Range("W3").Value = 1.99
$ws.Range("X3").Value = 23

# Row 6
$ws.Range("J6").Value = 3.5

# Row 7
$ws.Range("N7").Value = 5.4
$ws.Range("O7").Value = 1.21
$ws.Range("U7").Value = 2.2
$ws.Range("AJ7").Value = 14.5
$ws.Range("AM7").Value = 85

# Row 9
$ws.Range("Q9").Value = 1.92
$ws.Range("R9").Value = 1.41
$ws.Range("S9").Value = 3.3
$ws.Range("T9").Value = 1.75
$ws.Range("U9").Value = 2.26
$ws.Range("AA9").Value = 55
$ws.Range("AO9").Value = 32

# Row 10
$ws.Range("T10").Value = 1.91
$ws.Range("AF10").Value = 75
$ws.Range("AH10").Value = 23
